$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K6").Value = 0.5403464745801891

$ws.Range("J7").Value = 0.4969491838668565
$ws.Range("K7").Value = 0.2970525035592049

$ws.Range("I8").Value = 0.5011245128056051
$ws.Range("J8").Value = 0.2858677898194339

$ws.Range("H9").Value = 0.4852787037784192
$ws.Range("I9").Value = 0.2775335613519331

$ws.Range("G10").Value = 0.4539510573947921
$ws.Range("H10").Value = 0.2743085116504074

$ws.Range("F11").Value = 0.4663391832225094
$ws.Range("G11").Value = 0.2534447081011285

$ws.Range("E12").Value = 0.4814444548743619
$ws.Range("F12").Value = 0.2766837437271186

$ws.Range("D13").Value = 0.4184715358843989
$ws.Range("E13").Value = 0.2867219094086165

$ws.Range("C14").Value = 0.5177895860664353
$ws.Range("D14").Value = 0.1751453671933744

$ws.Range("B15").Value = 0.5618492773058843
$ws.Range("C15").Value = 0.1965658720679752

$ws.Range("B16").Value = 0.4328090033804217
